# Adiciona as colunas AR, AS e AT (31/12/2023, 31/03/2024, 30/06/2024)
# concatenando os balancos mais recentes na planilha existente (A1:AQ80 -> A1:AT80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cabecalho (linha 1): novas datas, com o mesmo estilo do cabecalho existente ---
$headerStyle = $ws.Range("AQ1").Style

$ws.Range("AR1").Value2 = "31/12/2023"
$ws.Range("AS1").Value2 = "31/03/2024"
$ws.Range("AT1").Value2 = "30/06/2024"

$headerRange = $ws.Range("AR1:AT1")
$headerRange.Style = $headerStyle
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin box)

# --- Linhas de dados (2-80): valores numericos nas 3 novas colunas ---
$ws.Range("AR2").Value2 = 2740000; $ws.Range("AS2").Value2 = 2658522.112; $ws.Range("AT2").Value2 = 2680698.112
$ws.Range("AR3").Value2 = 504263.008; $ws.Range("AS3").Value2 = 484182.016; $ws.Range("AT3").Value2 = 473985.984
$ws.Range("AR4").Value2 = 229159.008; $ws.Range("AS4").Value2 = 198348; $ws.Range("AT4").Value2 = 176167.008
$ws.Range("AR5").Value2 = 4950; $ws.Range("AS5").Value2 = 5109; $ws.Range("AT5").Value2 = 12720
$ws.Range("AR6").Value2 = 107685; $ws.Range("AS6").Value2 = 102956; $ws.Range("AT6").Value2 = 106271
$ws.Range("AR7").Value2 = 52765; $ws.Range("AS7").Value2 = 50902; $ws.Range("AT7").Value2 = 53246
$ws.Range("AR8").Value2 = 0; $ws.Range("AS8").Value2 = 0; $ws.Range("AT8").Value2 = 0
$ws.Range("AR9").Value2 = 63217; $ws.Range("AS9").Value2 = 69798; $ws.Range("AT9").Value2 = 62213
$ws.Range("AR10").Value2 = 9217; $ws.Range("AS10").Value2 = 14993; $ws.Range("AT10").Value2 = 19086
$ws.Range("AR11").Value2 = 37270; $ws.Range("AS11").Value2 = 42076; $ws.Range("AT11").Value2 = 44283
$ws.Range("AR12").Value2 = 135974; $ws.Range("AS12").Value2 = 139016.992; $ws.Range("AT12").Value2 = 127020
$ws.Range("AR13").Value2 = 6760; $ws.Range("AS13").Value2 = 6906; $ws.Range("AT13").Value2 = 0
$ws.Range("AR14").Value2 = 0; $ws.Range("AS14").Value2 = 0; $ws.Range("AT14").Value2 = 0
$ws.Range("AR15").Value2 = 0; $ws.Range("AS15").Value2 = 0; $ws.Range("AT15").Value2 = 0
$ws.Range("AR16").Value2 = 0; $ws.Range("AS16").Value2 = 0; $ws.Range("AT16").Value2 = 0
$ws.Range("AR17").Value2 = 0; $ws.Range("AS17").Value2 = 0; $ws.Range("AT17").Value2 = 0
$ws.Range("AR18").Value2 = 0; $ws.Range("AS18").Value2 = 0; $ws.Range("AT18").Value2 = 0
$ws.Range("AR19").Value2 = 49827; $ws.Range("AS19").Value2 = 59151; $ws.Range("AT19").Value2 = 66029
$ws.Range("AR20").Value2 = 0; $ws.Range("AS20").Value2 = 0; $ws.Range("AT20").Value2 = 0
$ws.Range("AR21").Value2 = 12026; $ws.Range("AS21").Value2 = 12458; $ws.Range("AT21").Value2 = 10076
$ws.Range("AR22").Value2 = 26585; $ws.Range("AS22").Value2 = 28556; $ws.Range("AT22").Value2 = 30404
$ws.Range("AR23").Value2 = 1120606.976; $ws.Range("AS23").Value2 = 1096816; $ws.Range("AT23").Value2 = 1120947.968
$ws.Range("AR24").Value2 = 952571.008; $ws.Range("AS24").Value2 = 909950.976; $ws.Range("AT24").Value2 = 928339.968
$ws.Range("AR25").Value2 = 0; $ws.Range("AS25").Value2 = 0; $ws.Range("AT25").Value2 = 0
$ws.Range("AR26").Value2 = 2740000; $ws.Range("AS26").Value2 = 2658522.112; $ws.Range("AT26").Value2 = 2680698.112
$ws.Range("AR27").Value2 = 547024; $ws.Range("AS27").Value2 = 466372.992; $ws.Range("AT27").Value2 = 480745.984
$ws.Range("AR28").Value2 = 83987; $ws.Range("AS28").Value2 = 80571; $ws.Range("AT28").Value2 = 65939
$ws.Range("AR29").Value2 = 239268; $ws.Range("AS29").Value2 = 167128; $ws.Range("AT29").Value2 = 175240
$ws.Range("AR30").Value2 = 12546; $ws.Range("AS30").Value2 = 4695; $ws.Range("AT30").Value2 = 7787
$ws.Range("AR31").Value2 = 76883; $ws.Range("AS31").Value2 = 76656; $ws.Range("AT31").Value2 = 72432
$ws.Range("AR32").Value2 = 0; $ws.Range("AS32").Value2 = 0; $ws.Range("AT32").Value2 = 0
$ws.Range("AR33").Value2 = 0; $ws.Range("AS33").Value2 = 0; $ws.Range("AT33").Value2 = 0
$ws.Range("AR34").Value2 = 134340; $ws.Range("AS34").Value2 = 137323.008; $ws.Range("AT34").Value2 = 159348
$ws.Range("AR35").Value2 = 0; $ws.Range("AS35").Value2 = 0; $ws.Range("AT35").Value2 = 0
$ws.Range("AR36").Value2 = 0; $ws.Range("AS36").Value2 = 0; $ws.Range("AT36").Value2 = 0
$ws.Range("AR37").Value2 = 1170125.056; $ws.Range("AS37").Value2 = 1185984; $ws.Range("AT37").Value2 = 1143042.048
$ws.Range("AR38").Value2 = 452867.008; $ws.Range("AS38").Value2 = 453375.008; $ws.Range("AT38").Value2 = 454108
$ws.Range("AR39").Value2 = 0; $ws.Range("AS39").Value2 = 0; $ws.Range("AT39").Value2 = 0
$ws.Range("AR40").Value2 = 578432; $ws.Range("AS40").Value2 = 564283.008; $ws.Range("AT40").Value2 = 581828.992
$ws.Range("AR41").Value2 = 40107; $ws.Range("AS41").Value2 = 40229; $ws.Range("AT41").Value2 = 43343
$ws.Range("AR42").Value2 = 0; $ws.Range("AS42").Value2 = 0; $ws.Range("AT42").Value2 = 0
$ws.Range("AR43").Value2 = 94662; $ws.Range("AS43").Value2 = 117241; $ws.Range("AT43").Value2 = 53372
$ws.Range("AR44").Value2 = 0; $ws.Range("AS44").Value2 = 0; $ws.Range("AT44").Value2 = 0
$ws.Range("AR45").Value2 = 4057; $ws.Range("AS45").Value2 = 10856; $ws.Range("AT45").Value2 = 10390
$ws.Range("AR46").Value2 = 0; $ws.Range("AS46").Value2 = 0; $ws.Range("AT46").Value2 = 0
$ws.Range("AR47").Value2 = 1022851.008; $ws.Range("AS47").Value2 = 1006164.992; $ws.Range("AT47").Value2 = 1056910.016
$ws.Range("AR48").Value2 = 1154461.952; $ws.Range("AS48").Value2 = 1154461.952; $ws.Range("AT48").Value2 = 1154461.952
$ws.Range("AR49").Value2 = 390692.992; $ws.Range("AS49").Value2 = 391755.008; $ws.Range("AT49").Value2 = 392807.008
$ws.Range("AR50").Value2 = 0; $ws.Range("AS50").Value2 = 0; $ws.Range("AT50").Value2 = 0
$ws.Range("AR51").Value2 = -5551; $ws.Range("AS51").Value2 = -5551; $ws.Range("AT51").Value2 = -5551
$ws.Range("AR52").Value2 = -570382.976; $ws.Range("AS52").Value2 = -598262.0159999999; $ws.Range("AT52").Value2 = -586064
$ws.Range("AR53").Value2 = 0; $ws.Range("AS53").Value2 = 0; $ws.Range("AT53").Value2 = 0
$ws.Range("AR54").Value2 = 0; $ws.Range("AS54").Value2 = 0; $ws.Range("AT54").Value2 = 0
$ws.Range("AR55").Value2 = 53630; $ws.Range("AS55").Value2 = 63761; $ws.Range("AT55").Value2 = 101256
$ws.Range("AR56").Value2 = 0; $ws.Range("AS56").Value2 = 0; $ws.Range("AT56").Value2 = 0

# Linha 57: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR57").Style = "Normal"
$ws.Range("AS57").Style = "Normal"
$ws.Range("AT57").Style = "Normal"

# Linha 58: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR58").Style = "Normal"
$ws.Range("AS58").Style = "Normal"
$ws.Range("AT58").Style = "Normal"
$ws.Range("AR59").Value2 = 464497.952; $ws.Range("AS59").Value2 = 506699.008; $ws.Range("AT59").Value2 = 570499.968
$ws.Range("AR60").Value2 = -300097.952; $ws.Range("AS60").Value2 = -361771.008; $ws.Range("AT60").Value2 = -379865.984
$ws.Range("AR61").Value2 = 164400; $ws.Range("AS61").Value2 = 144928; $ws.Range("AT61").Value2 = 190634
$ws.Range("AR62").Value2 = -175679.008; $ws.Range("AS62").Value2 = -142230; $ws.Range("AT62").Value2 = -161708.992
$ws.Range("AR63").Value2 = -62576.008; $ws.Range("AS63").Value2 = -46655; $ws.Range("AT63").Value2 = -58145
$ws.Range("AR64").Value2 = 0; $ws.Range("AS64").Value2 = 0; $ws.Range("AT64").Value2 = 0
$ws.Range("AR65").Value2 = 132530; $ws.Range("AS65").Value2 = 24943; $ws.Range("AT65").Value2 = 2704
$ws.Range("AR66").Value2 = -65125; $ws.Range("AS66").Value2 = 0; $ws.Range("AT66").Value2 = 0
$ws.Range("AR67").Value2 = 622; $ws.Range("AS67").Value2 = 2557; $ws.Range("AT67").Value2 = 60378
$ws.Range("AR68").Value2 = -44191; $ws.Range("AS68").Value2 = -23637; $ws.Range("AT68").Value2 = -21612
$ws.Range("AR69").Value2 = 0; $ws.Range("AS69").Value2 = 0; $ws.Range("AT69").Value2 = 0
$ws.Range("AR70").Value2 = 25953; $ws.Range("AS70").Value2 = 0; $ws.Range("AT70").Value2 = -21612

# Linha 71: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR71").Style = "Normal"
$ws.Range("AS71").Style = "Normal"
$ws.Range("AT71").Style = "Normal"

# Linha 72: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR72").Style = "Normal"
$ws.Range("AS72").Style = "Normal"
$ws.Range("AT72").Style = "Normal"

# Linha 73: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR73").Style = "Normal"
$ws.Range("AS73").Style = "Normal"
$ws.Range("AT73").Style = "Normal"
$ws.Range("AR74").Value2 = -50019; $ws.Range("AS74").Value2 = -40094; $ws.Range("AT74").Value2 = 12250
$ws.Range("AR75").Value2 = 13307; $ws.Range("AS75").Value2 = 3116; $ws.Range("AT75").Value2 = -3445
$ws.Range("AR76").Value2 = -50229; $ws.Range("AS76").Value2 = 9099; $ws.Range("AT76").Value2 = 3393

# Linha 77: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR77").Style = "Normal"
$ws.Range("AS77").Style = "Normal"
$ws.Range("AT77").Style = "Normal"

# Linha 78: linha em branco (rotulo de secao) no layout original -> mantem celulas vazias
$ws.Range("AR78").Style = "Normal"
$ws.Range("AS78").Style = "Normal"
$ws.Range("AT78").Style = "Normal"
$ws.Range("AR79").Value2 = 0; $ws.Range("AS79").Value2 = 0; $ws.Range("AT79").Value2 = 0
$ws.Range("AR80").Value2 = -76207; $ws.Range("AS80").Value2 = -27879; $ws.Range("AT80").Value2 = 12198
